{"js": "// Update racial-misclassification language to neutral \"50M voters\" phrasing.\n//\n// Three spots change (see XML diff):\n//   1. Professional-summary paragraph: \"... affecting all Black and\n//      Asian-American voters, developed ...\" -> \"... affecting 50M voters,\n//      developed ...\", single run, text-only change.\n//   2. Siege Analytics bullet: \"... affecting all Black and Asian-American\n//      voters, developed ...\" -> the old single run is split into three\n//      runs so the new \"50M\" sits in its own bold, dark-slate (#2C3E50) run,\n//      matching the bold \"23%\"/\"64%\" stat runs already in that bullet.\n//   3. \"Geospatial Demographic Classification System\" project Impact line:\n//      \"... affecting all Black and Asian-American voters, improved ...\" ->\n//      \"... affecting 50M voters nationwide, improved ...\", single run,\n//      text-only change (insertion of \" nationwide\" plus the phrase swap).\n//\n// All three paragraphs share the literal \"all Black and Asian-American\";\n// search for it (matchCase so we don't snag unrelated text) and handle each\n// hit according to which paragraph it falls in, then do a second, separate\n// search for the \" nationwide\" insertion in the Impact line.\n\nconst body = context.document.body;\n\nconst hits = body.search(\"all Black and Asian-American\", { matchCase: true });\nhits.load(\"text\");\nawait context.sync();\n\nfor (const hit of hits.items) {\n  const para = hit.paragraphs.getFirst();\n  para.load(\"text\");\n  await context.sync();\n\n  // Plain phrase swap everywhere; the bullet additionally gets bold +\n  // dark-slate color on just the new \"50M\" run.\n  hit.insertText(\"50M\", \"Replace\");\n\n  if (para.text.indexOf(\"race coding errors\") !== -1) {\n    hit.font.bold = true;\n    hit.font.color = \"#2C3E50\";\n  }\n}\n\nawait context.sync();\n\n// Impact line also gains \" nationwide\" right after \"voters,\" \u2014 done as a\n// follow-up pass once the phrase swap above has landed, searching on the\n// now-stable \"50M voters,\" text.\nconst impactHits = body.search(\"50M voters, improved electoral prediction accuracy\", {\n  matchCase: true,\n});\nawait context.sync();\n\nfor (const hit of impactHits.items) {\n  hit.insertText(\"50M voters nationwide, improved electoral prediction accuracy\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update racial-misclassification language to neutral \"50M voters\" phrasing.\n#\n# Three spots change (see XML diff):\n#   1. Professional-summary paragraph: \"... affecting all Black and\n#      Asian-American voters, developed ...\" -> \"... affecting 50M voters,\n#      developed ...\", single run, text-only change.\n#   2. Siege Analytics bullet: \"... affecting all Black and Asian-American\n#      voters, developed ...\" -> the old single run is split into three\n#      runs so the new \"50M\" sits in its own bold, dark-slate (#2C3E50) run,\n#      matching the bold \"23%\"/\"64%\" stat runs already in that bullet.\n#   3. \"Geospatial Demographic Classification System\" project Impact line:\n#      \"... affecting all Black and Asian-American voters, improved ...\" ->\n#      \"... affecting 50M voters nationwide, improved ...\", single run,\n#      text-only change (insertion of \" nationwide\" plus the phrase swap).\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n\n    if ($t -like \"*race coding errors affecting all Black and Asian-American voters*\") {\n        # Bullet #2 \u2014 isolate just the phrase, replace it, then bold + color\n        # that new run so it matches the other bold stat runs in the bullet.\n        $rng = $p.Range\n        $rng.Find.ClearFormatting()\n        $rng.Find.Replacement.ClearFormatting()\n        $found = $rng.Find.Execute(\"all Black and Asian-American\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n        if ($found) {\n            $rng.Text = \"50M\"\n            $rng.Font.Bold = 1\n            $rng.Font.Color = \"2C3E50\"\n        }\n    }\n    elseif ($t -like \"Impact: Corrected demographic data*\") {\n        # Project Impact line \u2014 phrase swap plus \" nationwide\" insertion.\n        $rng = $p.Range\n        $rng.Find.ClearFormatting()\n        $rng.Find.Replacement.ClearFormatting()\n        $rng.Find.Execute(\"all Black and Asian-American voters,\", $false, $false, $false, $false, $false, $true, 1, $false, \"50M voters nationwide,\", 2) | Out-Null\n    }\n    elseif ($t -like \"*affecting all Black and Asian-American voters*\") {\n        # Professional-summary sentence \u2014 plain phrase swap.\n        $rng = $p.Range\n        $rng.Find.ClearFormatting()\n        $rng.Find.Replacement.ClearFormatting()\n        $rng.Find.Execute(\"all Black and Asian-American\", $false, $false, $false, $false, $false, $true, 1, $false, \"50M\", 2) | Out-Null\n    }\n}\n"}
